# Add an "Avatar: {{avatar}}" block right after "Tuổi: {{age}}" and
# before "{{/user}}". The trailing "_GoBack" bookmark that used to sit
# at the end of the "Danh sách sinh viên:" paragraph is moved so that it
# now marks the insertion point of the new "{{avatar}}" text instead.

$d = $word.ActiveDocument

# Locate the "Tuổi: {{age}}" paragraph (4th paragraph in the document).
$agePara = $d.Paragraphs.Item(4)

# Insert a brand new paragraph right after it and fill it with "Avatar: ".
$insertionPoint = $agePara.Range
$insertionPoint.Collapse(0)   # wdCollapseEnd
$insertionPoint.InsertParagraphAfter()

$avatarLabelPara = $d.Paragraphs.Item(5)
$avatarLabelRange = $avatarLabelPara.Range
$avatarLabelRange.Collapse(0) # wdCollapseEnd
$avatarLabelRange.Text = "Avatar: "

# Insert another new paragraph after that one for the "{{avatar}}" tag.
$afterLabel = $avatarLabelPara.Range
$afterLabel.Collapse(0)       # wdCollapseEnd
$afterLabel.InsertParagraphAfter()

$avatarTagPara = $d.Paragraphs.Item(6)
$avatarTagRange = $avatarTagPara.Range
$avatarTagRange.Collapse(0)   # wdCollapseEnd
$avatarTagRange.Text = "{{avatar}}"

# Move the "_GoBack" bookmark from its old home (end of the "Danh sách
# sinh viên:" paragraph) to the start of the new "{{avatar}}" paragraph.
$d.Bookmarks.Item("_GoBack").Delete()

$avatarTagStart = $avatarTagPara.Range.Start
$bookmarkRange = $d.Range($avatarTagStart, $avatarTagStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
